$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 83000
$ws.Range("J48").Value = 83000
$ws.Range("L48").Value = 83000
$ws.Range("N48").Value = -83768
$ws.Range("H74").Value = 805.5476
$ws.Range("I74").Value = 722.2432
$ws.Range("J74").Value = 1422
$ws.Range("K74").Value = 722.2432
$ws.Range("L74").Value = 1422
$ws.Range("M74").Value = 151.7568
$ws.Range("N74").Value = -3170
$ws.Range("H77").Value = 805.5476
$ws.Range("I77").Value = 722.2432
$ws.Range("J77").Value = 1422
$ws.Range("K77").Value = 3611.216
$ws.Range("L77").Value = 7110
$ws.Range("M77").Value = 756.7840000000001
$ws.Range("N77").Value = -15846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1300
$ws.Range("I86").Value = 1350
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 1350
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -227
$ws.Range("N86").Value = -3446
$ws.Range("H89").Value = 1300
$ws.Range("I89").Value = 1350
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 6750
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -1134
$ws.Range("N89").Value = -17232
$ws.Range("H106").Value = 14250
$ws.Range("J106").Value = 14250
$ws.Range("L106").Value = 14250
$ws.Range("N106").Value = -16774
$ws.Range("H118").Value = 79000
$ws.Range("J118").Value = 79000
$ws.Range("L118").Value = 79000
$ws.Range("N118").Value = -82314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2860.1333
$ws.Range("I31").Value = 2449.5
$ws.Range("J31").Value = 3989.375
$ws.Range("K31").Value = 2449.5
$ws.Range("L31").Value = 3989.375
$ws.Range("M31").Value = -2154.5
$ws.Range("N31").Value = -4579.375
$ws.Range("H34").Value = 2860.1333
$ws.Range("I34").Value = 2449.5
$ws.Range("J34").Value = 3989.375
$ws.Range("K34").Value = 2449.5
$ws.Range("L34").Value = 3989.375
$ws.Range("M34").Value = -2247.5
$ws.Range("N34").Value = -4393.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 399.89474
$ws.Range("I5").Value = 367.7143
$ws.Range("J5").Value = 490
$ws.Range("K5").Value = 1103.1429
$ws.Range("L5").Value = 1470
$ws.Range("M5").Value = -991.1428999999998
$ws.Range("N5").Value = -1694
$ws.Range("H23").Value = 246
$ws.Range("I23").Value = 239.83333
$ws.Range("J23").Value = 250.625
$ws.Range("K23").Value = 719.49999
$ws.Range("L23").Value = 751.875
$ws.Range("M23").Value = -484.49999
$ws.Range("N23").Value = -1221.875
$ws.Range("H87").Value = 33860
$ws.Range("I87").Value = 1000
$ws.Range("J87").Value = 37511.11
$ws.Range("K87").Value = 3000
$ws.Range("L87").Value = 112533.33
$ws.Range("M87").Value = -1752
$ws.Range("N87").Value = -115029.33
$ws.Range("H90").Value = 33860
$ws.Range("I90").Value = 1000
$ws.Range("J90").Value = 37511.11
$ws.Range("K90").Value = 9000
$ws.Range("L90").Value = 337599.99
$ws.Range("M90").Value = -2760
$ws.Range("N90").Value = -350079.99
$ws.Range("H94").Value = 6286.4614
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 6286.4614
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = $null
$ws.Range("M94").Value = 18859.3842
$ws.Range("N94").Value = -20211.3842
$ws.Range("H95").Value = 6850.5713
$ws.Range("J95").Value = 6850.5713
$ws.Range("L95").Value = 20551.7139
$ws.Range("N95").Value = -24669.7139
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = $null
$ws.Range("N96").Value = 0
$ws.Range("H99").Value = 2001.9166
$ws.Range("I99").Value = 1346.1428
$ws.Range("J99").Value = 2920
$ws.Range("K99").Value = 4038.4284
$ws.Range("L99").Value = 8760
$ws.Range("M99").Value = -1792.4284
$ws.Range("N99").Value = -13252
$ws.Range("H100").Value = 6200
$ws.Range("J100").Value = 6200
$ws.Range("L100").Value = 18600
$ws.Range("N100").Value = -20222
$ws.Range("H103").Value = 2564.1667
$ws.Range("I103").Value = 461.66666
$ws.Range("J103").Value = 4666.6665
$ws.Range("K103").Value = 1384.99998
$ws.Range("L103").Value = 13999.9995
$ws.Range("M103").Value = -505.9999800000001
$ws.Range("N103").Value = -15757.9995
$ws.Range("H104").Value = 797.125
$ws.Range("I104").Value = 695
$ws.Range("J104").Value = 831.1667
$ws.Range("K104").Value = 2085
$ws.Range("L104").Value = 2493.5001
$ws.Range("M104").Value = 536
$ws.Range("N104").Value = -7735.5001
$ws.Range("H106").Value = 9000
$ws.Range("J106").Value = 9000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -28892
$ws.Range("H109").Value = 2152
$ws.Range("J109").Value = 5230
$ws.Range("L109").Value = 15690
$ws.Range("N109").Value = -17770
$ws.Range("H112").Value = 3608.5386
$ws.Range("J112").Value = 3608.5386
$ws.Range("L112").Value = 10825.6158
$ws.Range("N112").Value = -13041.6158
$ws.Range("H118").Value = 2372
$ws.Range("I118").Value = 2207
$ws.Range("K118").Value = 6621
$ws.Range("M118").Value = -5378
$ws.Range("H122").Value = 4950.32
$ws.Range("J122").Value = 6008.9
$ws.Range("L122").Value = 54080.1
$ws.Range("N122").Value = -58980.1
$ws.Range("H135").Value = 399.89474
$ws.Range("I135").Value = 367.7143
$ws.Range("J135").Value = 490
$ws.Range("K135").Value = 3309.4287
$ws.Range("L135").Value = 4410
$ws.Range("M135").Value = -774.4286999999999
$ws.Range("N135").Value = -9480
$ws.Range("H138").Value = 2048.7083
$ws.Range("I138").Value = 2027.5883
$ws.Range("J138").Value = 2100
$ws.Range("K138").Value = 6082.7649
$ws.Range("L138").Value = 6300
$ws.Range("M138").Value = -942.7649000000001
$ws.Range("N138").Value = -16580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 13500
$ws.Range("J27").Value = 13500
$ws.Range("L27").Value = 13500
$ws.Range("N27").Value = -13832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1099.7059
$ws.Range("J46").Value = 867.3333
$ws.Range("L46").Value = 867.3333
$ws.Range("N46").Value = -1243.3333
$ws.Range("H82").Value = 1907.5186
$ws.Range("I82").Value = 1300
$ws.Range("J82").Value = 2791.182
$ws.Range("K82").Value = 1300
$ws.Range("L82").Value = 2791.182
$ws.Range("M82").Value = -939
$ws.Range("N82").Value = -3513.182
$ws.Range("H85").Value = 1907.5186
$ws.Range("I85").Value = 1300
$ws.Range("J85").Value = 2791.182
$ws.Range("K85").Value = 1300
$ws.Range("L85").Value = 2791.182
$ws.Range("M85").Value = -52
$ws.Range("N85").Value = -5287.182
$ws.Range("H132").Value = 3140.1072
$ws.Range("I132").Value = 1877.6364
$ws.Range("J132").Value = 3957
$ws.Range("K132").Value = 5632.9092
$ws.Range("L132").Value = 11871
$ws.Range("M132").Value = -3102.9092
$ws.Range("N132").Value = -16931

Write-Host "Applied all changes"